$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23 (shifts existing rows 23-34 down to 24-35)
$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23, 1).Value = "C13"
$ws.Cells.Item(23, 2).Value = "C_0805_2012Metric"
$ws.Cells.Item(23, 3).Value = "100nF"
$ws.Cells.Item(23, 4).Value = "C28233"
$ws.Rows.Item(23).RowHeight = 13.5

$ws.Cells.Item(34, 4).Value = "C31850"

$ws.Range("D30").Select()
